# Apply data updates to the "Inscricoes" sheet as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: E=Inscritos, F=Pagos, G=Isenções deferidas, H=Inscrições homologadas

$ws.Cells.Item(2, 6).Value = 71
$ws.Cells.Item(2, 8).Value = 77

$ws.Cells.Item(7, 5).Value = 32

$ws.Cells.Item(10, 5).Value = 576
$ws.Cells.Item(10, 6).Value = 284
$ws.Cells.Item(10, 8).Value = 380

$ws.Cells.Item(11, 5).Value = 362

$ws.Cells.Item(12, 5).Value = 561
$ws.Cells.Item(12, 6).Value = 299
$ws.Cells.Item(12, 8).Value = 385

$ws.Cells.Item(15, 5).Value = 171
$ws.Cells.Item(15, 6).Value = 74
$ws.Cells.Item(15, 8).Value = 124

$ws.Cells.Item(16, 5).Value = 206

$ws.Cells.Item(21, 5).Value = 140

$ws.Cells.Item(22, 5).Value = 174
$ws.Cells.Item(22, 6).Value = 93
$ws.Cells.Item(22, 8).Value = 135

$ws.Cells.Item(23, 6).Value = 96
$ws.Cells.Item(23, 8).Value = 147

$ws.Cells.Item(24, 5).Value = 214

$ws.Cells.Item(25, 5).Value = 271
$ws.Cells.Item(25, 6).Value = 132
$ws.Cells.Item(25, 8).Value = 192

$ws.Cells.Item(26, 5).Value = 159

$ws.Cells.Item(27, 5).Value = 334

$ws.Cells.Item(29, 5).Value = 170
$ws.Cells.Item(29, 6).Value = 97
$ws.Cells.Item(29, 8).Value = 138

$ws.Cells.Item(30, 5).Value = 214

$ws.Cells.Item(33, 5).Value = 296
$ws.Cells.Item(33, 6).Value = 152
$ws.Cells.Item(33, 8).Value = 241

$ws.Cells.Item(35, 5).Value = 153

$ws.Cells.Item(36, 5).Value = 75

$ws.Cells.Item(37, 5).Value = 165
$ws.Cells.Item(37, 6).Value = 79
$ws.Cells.Item(37, 8).Value = 116

$ws.Cells.Item(38, 5).Value = 92

$ws.Cells.Item(40, 5).Value = 265

$ws.Cells.Item(41, 5).Value = 394

$ws.Cells.Item(42, 5).Value = 385

$ws.Cells.Item(44, 5).Value = 312
$ws.Cells.Item(44, 6).Value = 157
$ws.Cells.Item(44, 8).Value = 225

$ws.Cells.Item(46, 5).Value = 328
$ws.Cells.Item(46, 6).Value = 179
$ws.Cells.Item(46, 8).Value = 242

$ws.Cells.Item(47, 5).Value = 461
$ws.Cells.Item(47, 6).Value = 235
$ws.Cells.Item(47, 8).Value = 327

$ws.Cells.Item(50, 5).Value = 243
$ws.Cells.Item(50, 7).Value = 71
$ws.Cells.Item(50, 8).Value = 187

$ws.Cells.Item(51, 5).Value = 241
$ws.Cells.Item(51, 6).Value = 107
$ws.Cells.Item(51, 7).Value = 74
$ws.Cells.Item(51, 8).Value = 181

$ws.Cells.Item(52, 5).Value = 29
